# Applies the "add mrdodo and toringo scoring systems" edit:
#  - rows 2-5: mark C (Done?) as "yes" (first 4 tasks are done now); drop the old Queue (F) numbers
#  - append 6 new to-do rows (t5..t10) with their own Queue numbers
#  - widen column D, grow the table / autofilter / conditional formatting / data validation
#    ranges from A1:I5 (resp. C2:C5 / E2:E5) to A1:I11 (resp. C2:C11 / E2:E11)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Existing rows 2-5: these four tasks are now completed; clear the old
#    "Queue" position (column F) since finished items no longer queue up.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "yes"
$ws.Range("C3").Value = "yes"
$ws.Range("C4").Value = "yes"
$ws.Range("C5").Value = "yes"

$ws.Range("F2").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("F5").ClearContents()

# ---------------------------------------------------------------------------
# 2. New rows 6-11: the next batch of tennis-ranking tasks. Text is entered
#    in the same order the author originally typed it in (it was not a
#    strict top-to-bottom, left-to-right pass) so the shared-string table
#    comes out in the same order as the saved workbook.
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "t5"
$ws.Range("B6").Value = "Fai LogLikelihoodTerm.py"
$ws.Range("A7").Value = "t6"
$ws.Range("B7").Value = "Fai Loss.py"
$ws.Range("A8").Value = "t7"
$ws.Range("B8").Value = "Fai TennisDataframe.py"
$ws.Range("A9").Value = "t8"
$ws.Range("D8").Value = "Intendo la classe derivata da pd.Dataframe… con le colonne già inizializzate"
$ws.Range("A10").Value = "t9"
$ws.Range("B10").Value = "Fai import_notion_csv.py"
$ws.Range("D10").Value = "e genera il tennis dataframe"
$ws.Range("D9").Value = "in ingresso prende un tennis dataframe"
$ws.Range("B9").Value = "Fai TennisUniverse.py, escluso il metodo di ottimizzazione"
$ws.Range("A11").Value = "t10"
$ws.Range("B11").Value = "Fai il metodo di ottimizzazione di TennisUniverse"

$ws.Range("C6").Value = "no"
$ws.Range("C7").Value = "no"
$ws.Range("C8").Value = "no"
$ws.Range("C9").Value = "no"
$ws.Range("C10").Value = "no"
$ws.Range("C11").Value = "no"

$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 2
$ws.Range("F8").Value = 3
$ws.Range("F9").Value = 5
$ws.Range("F10").Value = 4
$ws.Range("F11").Value = 6

# ---------------------------------------------------------------------------
# 3. Grow the table (this also extends the AutoFilter range) to cover the
#    newly added rows.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:I11"))

# ---------------------------------------------------------------------------
# 4. Extend the conditional formatting on the "Done?" column down to the new
#    rows, keeping the same three rules (partially / yes / no).
# ---------------------------------------------------------------------------
$cfs = $ws.Range("C2:C5").FormatConditions
for ($i = 1; $i -le $cfs.Count; $i++) {
    $cfs.Item($i).ModifyAppliesToRange($ws.Range("C2:C11"))
}

# ---------------------------------------------------------------------------
# 5. Extend the data validation (dropdown) lists down to the new rows.
# ---------------------------------------------------------------------------
$ws.Range("C2:C5").Validation.Delete()
$ws.Range("C2:C11").Validation.Add(3, 1, 1, '"yes,partially,no"')

$ws.Range("E2:E5").Validation.Delete()
$ws.Range("E2:E11").Validation.Add(3, 1, 1, '"low;medium;high"')

# ---------------------------------------------------------------------------
# 6. Cosmetic follow-on: column D got a bit wider to fit the new notes, and
#    the active selection moved to where the user continued working.
# ---------------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 67.21875
$ws.Range("B22").Select()
